# Updates cryptos list (commit: "Updated cryptos list on Fri Nov  1 12:17:26 UTC 2024 with GitHub Actions")
# Refreshes Price (col D) and Volume(1h) (col E) figures for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.179.94"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "2.520.83"
$ws.Range("E3").Value = "  -4.36%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.76%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "2.519.21"
$ws.Range("E9").Value = "  -4.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "2.980.07"
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").Value = "70.015.27"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("E16").Value = "  -4.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "2.520.04"
$ws.Range("E18").Value = "  -4.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.67%  "
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("D28").Value = "2.647.99"
$ws.Range("E28").Value = "  -4.54%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "0.0₃0913"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "465.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.15%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "154.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.36%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.322"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("E43").Value = "  -6.74%  "
$ws.Range("E44").Value = "  -13.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.532"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("E50").Value = "  -4.08%  "
$ws.Range("E51").Value = "  -1.02%  "
